$wb = $excel.ActiveWorkbook

# --- Update values on the existing "Delivery_results" sheet ---
$ws1 = $wb.Worksheets.Item("Delivery_results")

$ws1.Range("F2").Value = 0.51
$ws1.Range("G2").Value = 4.81

$ws1.Range("D3").Value = 44
$ws1.Range("E3").Value = 73.33
$ws1.Range("F3").Value = 0.1
$ws1.Range("G3").Value = 1.41

$ws1.Range("F4").Value = 0.26
$ws1.Range("G4").Value = 2.35

$ws1.Range("F5").Value = 0.93
$ws1.Range("G5").Value = 6.22

$ws1.Range("F6").Value = 0.47
$ws1.Range("G6").Value = 4.38

# --- Add a new "Total_distance" sheet after "Total_emissions" ---
$wsEmissions = $wb.Worksheets.Item("Total_emissions")
$wsNew = $wb.Worksheets.Add($null, $wsEmissions)
$wsNew.Name = "Total_distance"

# Copy the header cell formatting (bold, centered, bordered) from the
# "Total_emissions" sheet's A1 header cell, then overwrite with our own values.
$wsEmissions.Range("A1").Copy($wsNew.Range("A1"))
$wsNew.Range("A1").Value = "Total distance (km)"
$wsNew.Range("A2").Value = 39.76204
